# Fix the display issue reported in the commit: update the first
# reservation row and append the missing reservation rows that were
# not being rendered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: replace the outdated "karoui / 2025-02-28" entry with the
#     actual first reservation (raslen). Status stays "Confirmée".
$ws.Range("A3").Value = "raslen"
$ws.Range("B3").Value = "2025-04-27 19:04:00"
$ws.Range("C3").Value = "Confirmée"

# --- Rows 4-18 are untouched (still ahmed / karoui history) ---

# --- New row 19: continues the existing "Confirmé" (accented, filled)
#     styled block, same as rows 6-18.
$ws.Range("A19").Value = "karoui"
$ws.Range("B19").Value = "2025-03-07 09:52:35"
$ws.Range("B19").HorizontalAlignment = -4108
$ws.Range("C19").Value = "Confirmé"
$ws.Range("C19").Interior.ColorIndex = 35

# --- New rows 20-25: newest reservations, status "confirmée" (lower
#     case) with plain/default styling (no fill).
$newRows = @(
    @{ Row = 20; Date = "2025-04-20 09:31:18" },
    @{ Row = 21; Date = "2025-04-20 09:35:25" },
    @{ Row = 22; Date = "2025-04-20 11:19:35" },
    @{ Row = 23; Date = "2025-04-20 11:21:09" },
    @{ Row = 24; Date = "2025-04-22 08:22:58" },
    @{ Row = 25; Date = "2025-04-22 18:50:50" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = "karoui"
    $ws.Range("B$row").Value = $r.Date
    $ws.Range("B$row").HorizontalAlignment = -4108
    $ws.Range("C$row").Value = "confirmée"
}
